$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update C2 time value ---
$ws.Cells.Item(2, 3).Value = "11:00:00"

# --- Update changed numeric cells in rows 2-14 ---
# Row 2
$ws.Cells.Item(2, 19).Value = 1.05
$ws.Cells.Item(2, 22).Value = 1.19
$ws.Cells.Item(2, 23).Value = 1.01

# Row 3
$ws.Cells.Item(3, 8).Value = 10.5
$ws.Cells.Item(3, 9).Value = 11.5
$ws.Cells.Item(3, 11).Value = 6
$ws.Cells.Item(3, 16).Value = 2.64
$ws.Cells.Item(3, 25).Value = 42
$ws.Cells.Item(3, 26).Value = 100
$ws.Cells.Item(3, 28).Value = 11
$ws.Cells.Item(3, 29).Value = 13
$ws.Cells.Item(3, 31).Value = 150
$ws.Cells.Item(3, 34).Value = 25
$ws.Cells.Item(3, 35).Value = 120
$ws.Cells.Item(3, 41).Value = 160

# Row 4
$ws.Cells.Item(4, 7).Value = 5.4
$ws.Cells.Item(4, 9).Value = 1.75
$ws.Cells.Item(4, 22).Value = 2.32
$ws.Cells.Item(4, 23).Value = 1.22
$ws.Cells.Item(4, 29).Value = 9.4
$ws.Cells.Item(4, 34).Value = 17.5
$ws.Cells.Item(4, 39).Value = 80
$ws.Cells.Item(4, 41).Value = 8.199999999999999

# Row 5
$ws.Cells.Item(5, 18).Value = 1.17

# Row 6
$ws.Cells.Item(6, 25).Value = 980
$ws.Cells.Item(6, 26).Value = 980
$ws.Cells.Item(6, 28).Value = 16
$ws.Cells.Item(6, 29).Value = 10.5
$ws.Cells.Item(6, 31).Value = 980
$ws.Cells.Item(6, 33).Value = 13
$ws.Cells.Item(6, 35).Value = 980
$ws.Cells.Item(6, 37).Value = 980
$ws.Cells.Item(6, 38).Value = 980
$ws.Cells.Item(6, 40).Value = 11
$ws.Cells.Item(6, 41).Value = 1000

# Row 7
$ws.Cells.Item(7, 6).Value = 2.78
$ws.Cells.Item(7, 7).Value = 2.8
$ws.Cells.Item(7, 9).Value = 2.66
$ws.Cells.Item(7, 22).Value = 1.6
$ws.Cells.Item(7, 23).Value = 1.55
$ws.Cells.Item(7, 27).Value = 38
$ws.Cells.Item(7, 28).Value = 15
$ws.Cells.Item(7, 36).Value = 40
$ws.Cells.Item(7, 39).Value = 60
$ws.Cells.Item(7, 40).Value = 17
$ws.Cells.Item(7, 41).Value = 16

# Row 8
$ws.Cells.Item(8, 7).Value = 8.199999999999999
$ws.Cells.Item(8, 10).Value = 5.5
$ws.Cells.Item(8, 16).Value = 2.76
$ws.Cells.Item(8, 20).Value = 1.79
$ws.Cells.Item(8, 23).Value = 1.14
$ws.Cells.Item(8, 27).Value = 13
$ws.Cells.Item(8, 37).Value = 95
$ws.Cells.Item(8, 38).Value = 80
$ws.Cells.Item(8, 39).Value = 95

# Row 9
$ws.Cells.Item(9, 6).Value = 2.28
$ws.Cells.Item(9, 7).Value = 2.3
$ws.Cells.Item(9, 8).Value = 3.35
$ws.Cells.Item(9, 12).Value = 1.35
$ws.Cells.Item(9, 13).Value = 1.06
$ws.Cells.Item(9, 17).Value = 1.8
$ws.Cells.Item(9, 18).Value = 1.47
$ws.Cells.Item(9, 21).Value = 2.42
$ws.Cells.Item(9, 23).Value = 1.76
$ws.Cells.Item(9, 24).Value = 16.5
$ws.Cells.Item(9, 26).Value = 24
$ws.Cells.Item(9, 28).Value = 12
$ws.Cells.Item(9, 30).Value = 13.5
$ws.Cells.Item(9, 34).Value = 15.5
$ws.Cells.Item(9, 36).Value = 29
$ws.Cells.Item(9, 39).Value = 75
$ws.Cells.Item(9, 41).Value = 28

# Row 10
$ws.Cells.Item(10, 12).Value = 1.29
$ws.Cells.Item(10, 15).Value = 1.2
$ws.Cells.Item(10, 17).Value = 1.6
$ws.Cells.Item(10, 21).Value = 2.46
$ws.Cells.Item(10, 28).Value = 12
$ws.Cells.Item(10, 29).Value = 10.5

# Row 11
$ws.Cells.Item(11, 7).Value = 1.4
$ws.Cells.Item(11, 9).Value = 10.5
$ws.Cells.Item(11, 10).Value = 5.5
$ws.Cells.Item(11, 11).Value = 5.6
$ws.Cells.Item(11, 15).Value = 1.17
$ws.Cells.Item(11, 16).Value = 2.74
$ws.Cells.Item(11, 17).Value = 1.55
$ws.Cells.Item(11, 19).Value = 2.36
$ws.Cells.Item(11, 20).Value = 1.85
$ws.Cells.Item(11, 23).Value = 3.5
$ws.Cells.Item(11, 25).Value = 40
$ws.Cells.Item(11, 40).Value = 4.6

# Row 12
$ws.Cells.Item(12, 6).Value = 1.33
$ws.Cells.Item(12, 9).Value = 10.5
$ws.Cells.Item(12, 11).Value = 6.6
$ws.Cells.Item(12, 21).Value = 2.46
$ws.Cells.Item(12, 23).Value = 3.85
$ws.Cells.Item(12, 39).Value = 75
$ws.Cells.Item(12, 40).Value = 3.45
$ws.Cells.Item(12, 41).Value = 70

# Row 13
$ws.Cells.Item(13, 7).Value = 2.4
$ws.Cells.Item(13, 9).Value = 3.3
$ws.Cells.Item(13, 17).Value = 1.9
$ws.Cells.Item(13, 18).Value = 1.42
$ws.Cells.Item(13, 23).Value = 1.71
$ws.Cells.Item(13, 34).Value = 16
$ws.Cells.Item(13, 40).Value = 17

# Row 14
$ws.Cells.Item(14, 6).Value = 2.02
$ws.Cells.Item(14, 8).Value = 3.7
$ws.Cells.Item(14, 9).Value = 4.1
$ws.Cells.Item(14, 14).Value = 4.1
$ws.Cells.Item(14, 18).Value = 1.42
$ws.Cells.Item(14, 19).Value = 3
$ws.Cells.Item(14, 20).Value = 1.71
$ws.Cells.Item(14, 21).Value = 2.2
$ws.Cells.Item(14, 22).Value = 1.32
$ws.Cells.Item(14, 23).Value = 1.87
$ws.Cells.Item(14, 24).Value = 21
$ws.Cells.Item(14, 25).Value = 16
$ws.Cells.Item(14, 26).Value = 30
$ws.Cells.Item(14, 27).Value = 75
$ws.Cells.Item(14, 28).Value = 11
$ws.Cells.Item(14, 29).Value = 8.800000000000001
$ws.Cells.Item(14, 30).Value = 16
$ws.Cells.Item(14, 31).Value = 44
$ws.Cells.Item(14, 32).Value = 14
$ws.Cells.Item(14, 33).Value = 11
$ws.Cells.Item(14, 34).Value = 17
$ws.Cells.Item(14, 35).Value = 50
$ws.Cells.Item(14, 36).Value = 25
$ws.Cells.Item(14, 37).Value = 21
$ws.Cells.Item(14, 38).Value = 34
$ws.Cells.Item(14, 39).Value = 100
$ws.Cells.Item(14, 40).Value = 14
$ws.Cells.Item(14, 41).Value = 42

# --- Insert new row 15 (Chilean Primera B) and shift old row 15 to row 16 ---
$ws.Rows.Item(15).Insert()

# --- Set values for new row 15 ---
$ws.Cells.Item(15, 1).Value = "Chilean Primera B"
$ws.Cells.Item(15, 2).Value = "'2025-11-26"
$ws.Cells.Item(15, 3).Value = "20:30:00"
$ws.Cells.Item(15, 4).Value = "Deportes Concepcion"
$ws.Cells.Item(15, 5).Value = "Deportes Copiapo"
$ws.Cells.Item(15, 6).Value = 1.09
$ws.Cells.Item(15, 7).Value = 1000
$ws.Cells.Item(15, 8).Value = 1.09
$ws.Cells.Item(15, 9).Value = 1000
$ws.Cells.Item(15, 10).Value = 1.09
$ws.Cells.Item(15, 11).Value = 1000
$ws.Cells.Item(15, 12).Value = 1.01
$ws.Cells.Item(15, 13).Value = 1.01
$ws.Cells.Item(15, 14).Value = 1.1
$ws.Cells.Item(15, 15).Value = 1.01
$ws.Cells.Item(15, 16).Value = 1.24
$ws.Cells.Item(15, 17).Value = 1.02
$ws.Cells.Item(15, 18).Value = 1.18
$ws.Cells.Item(15, 19).Value = 1.42
$ws.Cells.Item(15, 20).Value = 1.04
$ws.Cells.Item(15, 21).Value = 1.04
$ws.Cells.Item(15, 22).Value = 1.01
$ws.Cells.Item(15, 23).Value = 1.01
$ws.Cells.Item(15, 24).Value = 1000
$ws.Cells.Item(15, 25).Value = 1000
$ws.Cells.Item(15, 26).Value = 1000
$ws.Cells.Item(15, 27).Value = 1000
$ws.Cells.Item(15, 28).Value = 1000
$ws.Cells.Item(15, 29).Value = 1000
$ws.Cells.Item(15, 30).Value = 1000
$ws.Cells.Item(15, 31).Value = 1000
$ws.Cells.Item(15, 32).Value = 1000
$ws.Cells.Item(15, 33).Value = 1000
$ws.Cells.Item(15, 34).Value = 1000
$ws.Cells.Item(15, 35).Value = 1000
$ws.Cells.Item(15, 36).Value = 1000
$ws.Cells.Item(15, 37).Value = 1000
$ws.Cells.Item(15, 38).Value = 1000
$ws.Cells.Item(15, 39).Value = 1000
$ws.Cells.Item(15, 40).Value = 1000
$ws.Cells.Item(15, 41).Value = 1000

# --- Set values for row 16 (previously row 15, Colombian Primera A) ---
$ws.Cells.Item(16, 1).Value = "Colombian Primera A"
$ws.Cells.Item(16, 2).Value = "'2025-11-26"
$ws.Cells.Item(16, 3).Value = "20:30:00"
$ws.Cells.Item(16, 4).Value = "Atletico Bucaramanga"
$ws.Cells.Item(16, 5).Value = "Fortaleza FC"
$ws.Cells.Item(16, 6).Value = 1.09
$ws.Cells.Item(16, 7).Value = 2
$ws.Cells.Item(16, 8).Value = 4
$ws.Cells.Item(16, 9).Value = 7.8
$ws.Cells.Item(16, 10).Value = 1.09
$ws.Cells.Item(16, 11).Value = 980
$ws.Cells.Item(16, 12).Value = 1.01
$ws.Cells.Item(16, 13).Value = 1.01
$ws.Cells.Item(16, 14).Value = 2.44
$ws.Cells.Item(16, 15).Value = 1.02
$ws.Cells.Item(16, 16).Value = 1.54
$ws.Cells.Item(16, 17).Value = 2.04
$ws.Cells.Item(16, 18).Value = 1.19
$ws.Cells.Item(16, 19).Value = 3.55
$ws.Cells.Item(16, 20).Value = 1.04
$ws.Cells.Item(16, 21).Value = 1.04
$ws.Cells.Item(16, 22).Value = 1.18
$ws.Cells.Item(16, 23).Value = 2
$ws.Cells.Item(16, 24).Value = 1000
$ws.Cells.Item(16, 25).Value = 1000
$ws.Cells.Item(16, 26).Value = 1000
$ws.Cells.Item(16, 27).Value = 1000
$ws.Cells.Item(16, 28).Value = 1000
$ws.Cells.Item(16, 29).Value = 1000
$ws.Cells.Item(16, 30).Value = 1000
$ws.Cells.Item(16, 31).Value = 1000
$ws.Cells.Item(16, 32).Value = 1000
$ws.Cells.Item(16, 33).Value = 1000
$ws.Cells.Item(16, 34).Value = 1000
$ws.Cells.Item(16, 35).Value = 1000
$ws.Cells.Item(16, 36).Value = 1000
$ws.Cells.Item(16, 37).Value = 1000
$ws.Cells.Item(16, 38).Value = 1000
$ws.Cells.Item(16, 39).Value = 1000
$ws.Cells.Item(16, 40).Value = 1000
$ws.Cells.Item(16, 41).Value = 1000